$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rf")

$ws.Range("B2").Value = 0.7844141069397043

$report = "              precision    recall  f1-score   support`n`n           0       0.79      0.96      0.87      1300`n           1       0.73      0.28      0.40       458`n`n    accuracy                           0.78      1758`n   macro avg       0.76      0.62      0.63      1758`nweighted avg       0.77      0.78      0.75      1758`n"
$ws.Range("C2").Value = $report

$cm = "[[1253   47]`n [ 332  126]]"
$ws.Range("D2").Value = $cm

$ws.Rows.Item(2).AutoFit()
